# Automated MeteoCat daily-summary refresh.
# Re-applies the latest extraction pass (new DATA_EXTRACCIO timestamps)
# together with the handful of metric values that moved between the
# previous run and this one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Row 2
$ws.Range("E2").Value = "2026-02-16 21:48:53"
$ws.Range("H2").Value = "'95%"
$ws.Range("N2").Value = "0.6 °C 21:29 TU"
# Row 3
$ws.Range("E3").Value = "2026-02-16 21:48:56"
$ws.Range("L3").Value = "72.7 km/h - 233º 21:02 TU"
$ws.Range("N3").Value = "-4.1 °C 21:29 TU"
# Row 4
$ws.Range("E4").Value = "2026-02-16 21:48:59"
$ws.Range("H4").Value = "'62%"
$ws.Range("N4").Value = "7.0 °C 21:25 TU"
$ws.Range("O4").Value = "13.4 °C"
# Row 5
$ws.Range("E5").Value = "2026-02-16 21:49:02"
$ws.Range("N5").Value = "-3.8 °C 21:29 TU"
$ws.Range("O5").Value = "-0.8 °C"
# Row 6
$ws.Range("E6").Value = "2026-02-16 21:49:04"
# Row 7
$ws.Range("E7").Value = "2026-02-16 21:49:07"
# Row 8
$ws.Range("E8").Value = "2026-02-16 21:49:10"
# Row 9
$ws.Range("E9").Value = "2026-02-16 21:49:12"
$ws.Range("L9").Value = "50.0 km/h - 352º 21:09 TU"
# Row 10
$ws.Range("E10").Value = "2026-02-16 21:49:15"
# Row 11
$ws.Range("E11").Value = "2026-02-16 21:49:17"
# Row 12
$ws.Range("E12").Value = "2026-02-16 21:49:20"
$ws.Range("H12").Value = "'78%"
$ws.Range("O12").Value = "11.1 °C"
# Row 13
$ws.Range("E13").Value = "2026-02-16 21:49:23"
# Row 14
$ws.Range("E14").Value = "2026-02-16 21:49:25"
# Row 15
$ws.Range("E15").Value = "2026-02-16 21:49:28"
# Row 16
$ws.Range("E16").Value = "2026-02-16 21:49:31"
$ws.Range("N16").Value = "-3.4 °C 21:28 TU"
$ws.Range("O16").Value = "-0.3 °C"
# Row 17
$ws.Range("E17").Value = "2026-02-16 21:49:34"
# Row 18
$ws.Range("E18").Value = "2026-02-16 21:49:37"
$ws.Range("O18").Value = "10.7 °C"
# Row 19
$ws.Range("E19").Value = "2026-02-16 21:49:39"
# Row 20
$ws.Range("E20").Value = "2026-02-16 21:49:42"
$ws.Range("L20").Value = "61.6 km/h - 331º 21:18 TU"
$ws.Range("N20").Value = "-2.1 °C 21:29 TU"
# Row 21
$ws.Range("E21").Value = "2026-02-16 21:49:45"
# Row 22
$ws.Range("E22").Value = "2026-02-16 21:49:48"
# Row 23
$ws.Range("E23").Value = "2026-02-16 21:49:51"
$ws.Range("I23").Value = "16.1 mm"
$ws.Range("N23").Value = "-4.6 °C 21:29 TU"
$ws.Range("O23").Value = "-0.8 °C"
# Row 24
$ws.Range("E24").Value = "2026-02-16 21:49:53"
$ws.Range("H24").Value = "'72%"
$ws.Range("J24").Value = "1016.4 hPa"
# Row 25
$ws.Range("E25").Value = "2026-02-16 21:49:56"
$ws.Range("I25").Value = "6.6 mm"
$ws.Range("N25").Value = "-1.3 °C 21:29 TU"
# Row 26
$ws.Range("E26").Value = "2026-02-16 21:49:59"
# Row 27
$ws.Range("E27").Value = "2026-02-16 21:50:02"
$ws.Range("N27").Value = "-0.2 °C 21:25 TU"
# Row 28
$ws.Range("E28").Value = "2026-02-16 21:50:04"
# Row 29
$ws.Range("E29").Value = "2026-02-16 21:50:07"
$ws.Range("O29").Value = "10.8 °C"
# Row 30
$ws.Range("E30").Value = "2026-02-16 21:50:10"
$ws.Range("H30").Value = "'68%"
$ws.Range("L30").Value = "43.6 km/h - 30º 21:11 TU"
$ws.Range("O30").Value = "11.9 °C"
# Row 31
$ws.Range("E31").Value = "2026-02-16 21:50:13"
$ws.Range("N31").Value = "11.4 °C 21:14 TU"
$ws.Range("O31").Value = "14.4 °C"
# Row 32
$ws.Range("E32").Value = "2026-02-16 21:50:15"
# Row 33
$ws.Range("E33").Value = "2026-02-16 21:50:18"
$ws.Range("L33").Value = "33.8 km/h - 324º 21:28 TU"
# Row 34
$ws.Range("E34").Value = "2026-02-16 21:50:21"
$ws.Range("N34").Value = "1.3 °C 21:29 TU"
$ws.Range("O34").Value = "3.4 °C"
# Row 35
$ws.Range("E35").Value = "2026-02-16 21:50:23"
$ws.Range("I35").Value = "2.2 mm"
# Row 36
$ws.Range("E36").Value = "2026-02-16 21:50:26"
$ws.Range("H36").Value = "'71%"
$ws.Range("O36").Value = "12.2 °C"
# Row 37
$ws.Range("E37").Value = "2026-02-16 21:50:28"
# Row 38
$ws.Range("E38").Value = "2026-02-16 21:50:31"
# Row 39
$ws.Range("E39").Value = "2026-02-16 21:50:34"
$ws.Range("I39").Value = "4.7 mm"
$ws.Range("N39").Value = "-3.5 °C 21:29 TU"
$ws.Range("O39").Value = "0.1 °C"
# Row 40
$ws.Range("E40").Value = "2026-02-16 21:50:37"
# Row 41
$ws.Range("E41").Value = "2026-02-16 21:50:39"
$ws.Range("O41").Value = "17.7 °C"
# Row 42
$ws.Range("E42").Value = "2026-02-16 21:50:42"
$ws.Range("I42").Value = "0.1 mm"
$ws.Range("O42").Value = "11.4 °C"
# Row 43
$ws.Range("E43").Value = "2026-02-16 21:50:44"
# Row 44
$ws.Range("E44").Value = "2026-02-16 21:50:46"
$ws.Range("G44").Value = "244 cm"
$ws.Range("H44").Value = "'90%"
$ws.Range("I44").Value = "13.4 mm"
$ws.Range("N44").Value = "-3.2 °C 21:29 TU"
$ws.Range("O44").Value = "-0.3 °C"
# Row 45
$ws.Range("E45").Value = "2026-02-16 21:50:49"
$ws.Range("L45").Value = "31.0 km/h - 230º 21:13 TU"
$ws.Range("O45").Value = "4.8 °C"
# Row 46
$ws.Range("E46").Value = "2026-02-16 21:50:52"
